$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the style of the
# existing header row (bold font, border, centered alignment - same as H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data values for columns I and J, rows 2-4
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8

$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 9

$ws.Range("I4").Value = 9
$ws.Range("J4").Value = 9
